$d = $word.ActiveDocument

# 1. Update the date on the cover page
$d.Content.Find.Execute("Date: 2024-09-11", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Date: 2024-09-12", 2)

# 2. Table 1: "Overview of the nodes in the control flow"
$t1 = $d.Tables.Item(1)

# Execute SQL Task occurrences: 1 -> 4
$t1.Cell(2, 2).Range.Text = "4"

# Insert a new "Data Flow Task" row before the "Foreach Loop Container" row
$beforeRow = $t1.Rows.Item(3)
$newRow = $t1.Rows.Add($beforeRow)
$newRow.Cells.Item(1).Range.Text = "Data Flow Task"
$newRow.Cells.Item(2).Range.Text = "1"

# "Foreach Loop Container" -> "Expression Task" (now row 4 after the insert)
$t1.Cell(4, 1).Range.Text = "Expression Task"

# 3. Table 3: "Overview of the nodes in the data flow" - add new rows after the header
$t3 = $d.Tables.Item(3)
$dataFlowNodes = @(
    @("DataSources", "3"),
    @("Variable", "2"),
    @("DerivedColumn", "2"),
    @("RowCount", "2"),
    @("SSISODBCDst", "2"),
    @("DataDestinations", "1"),
    @("SSISODBCSrc", "1"),
    @("Lookup", "1"),
    @("ConditionalSplit", "1"),
    @("UnionAll", "1")
)
foreach ($pair in $dataFlowNodes) {
    $row = $t3.Rows.Add()
    $row.Cells.Item(1).Range.Text = $pair[0]
    $row.Cells.Item(2).Range.Text = $pair[1]
}

# 4. Table 4: "Overview of utilised source tables in the data flow" - add new rows
$t4 = $d.Tables.Item(4)
$sourceTables = @(
    @("Suppliers_Extract", "1"),
    @("Products_Extract", "1")
)
foreach ($pair in $sourceTables) {
    $row = $t4.Rows.Add()
    $row.Cells.Item(1).Range.Text = $pair[0]
    $row.Cells.Item(2).Range.Text = $pair[1]
}

# 5. Table 5: "Overview of utilised target tables in the data flow" - add new row
$t5 = $d.Tables.Item(5)
$row = $t5.Rows.Add()
$row.Cells.Item(1).Range.Text = "Error_lines"
$row.Cells.Item(2).Range.Text = "1"

# 6. Table 6: transformations table - replace "nan" placeholders with real node names
$t6 = $d.Tables.Item(6)
$t6.Cell(2, 1).Range.Text = "Merge and filter@Alter name pref"
$t6.Cell(3, 1).Range.Text = "Merge and filter@Error_match_column"

# 7. Add a blank paragraph after the "Sankey Diagrams" heading (before the
#    description paragraph), matching the plain, unstyled paragraphs used
#    elsewhere in this document.
$bodyParagraphs = $d.Content.Paragraphs
for ($i = 1; $i -le $bodyParagraphs.Count; $i++) {
    $para = $bodyParagraphs.Item($i)
    if ($para.Range.Text -match "^Sankey Diagrams") {
        $nextPara = $bodyParagraphs.Item($i + 1)
        $nextPara.Range.InsertParagraphBefore()
        break
    }
}

# 8. Fix casing: "a Sankey Diagram" -> "a sankey Diagram" in the description paragraph
$d.Content.Find.Execute("This section contains the Merge and filter data flow in a Sankey Diagram,", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This section contains the Merge and filter data flow in a sankey Diagram,", 2)

# 9. Legend table: "- join or split node" -> "- Join or split node"
$d.Content.Find.Execute("- join or split node", $true, $false, $false, $false, $false,
                         $true, 1, $false, "- Join or split node", 2)
